$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ================= Header text: Volume Number 7 -> 8 =================
$hdr1 = $ws.Range("A8")
$hdr1Whole = $hdr1.Characters(1, 21)
$hdr1Whole.Text = "Volume 31   Number  8"
foreach ($seg in @(@(1,7), @(8,2), @(10,11), @(21,1))) {
    $run = $hdr1.Characters($seg[0], $seg[1])
    $run.Font.Name = "Andale WT"
    $run.Font.Size = 10
}

# ========= Header text: Report week 2/12-2/18/2024 -> 2/19-2/25/2024 =========
$hdr2 = $ws.Range("C9")
$hdr2Whole = $hdr2.Characters(1, 55)
$hdr2Whole.Text = "Report Covering the Week  2/19/2024  Through  2/25/2024"
foreach ($seg in @(@(1,26), @(27,9), @(36,11), @(47,9))) {
    $run = $hdr2.Characters($seg[0], $seg[1])
    $run.Font.Name = "Andale WT"
    $run.Font.Size = 10
}

# ================= Weekly crime-stat table updates (rows 14-30) =================

# ---- Cells changing from placeholder text "0" (style General) to a real count (style #,##0) ----
$ws.Range("F15").Copy() | Out-Null
foreach ($ref in @("C15", "C16", "D22", "G22", "C26", "C27", "D27", "C30", "F30", "I30")) {
    $ws.Range($ref).PasteSpecial(-4122) | Out-Null
}
$ws.Range("C15").Value = 1
$ws.Range("C16").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("G22").Value = 1
$ws.Range("C26").Value = 1
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("C30").Value = 1
$ws.Range("F30").Value = 1
$ws.Range("I30").Value = 1

# ---- Cells changing from a real count (style #,##0) back to placeholder text "0" ----
$ws.Range("C14").Copy() | Out-Null
foreach ($ref in @("C17", "D20")) {
    $ws.Range($ref).PasteSpecial(-4122) | Out-Null
    $ws.Range("C14").Copy() | Out-Null
    $ws.Range($ref).PasteSpecial(-4163) | Out-Null
}

# ---- Cells changing from placeholder text "***.*" (style General) to a real pct (style #,##0.0) ----
$ws.Range("H15").Copy() | Out-Null
foreach ($ref in @("M14", "E22", "H22", "E27")) {
    $ws.Range($ref).PasteSpecial(-4122) | Out-Null
}
$ws.Range("M14").Value = -100
$ws.Range("E22").Value = -100
$ws.Range("H22").Value = 0
$ws.Range("E27").Value = 0

# ---- Cells changing from a real pct (style #,##0.0) back to placeholder text "***.*" ----
$ws.Range("E14").Copy() | Out-Null
foreach ($ref in @("E20")) {
    $ws.Range($ref).PasteSpecial(-4122) | Out-Null
    $ws.Range("E14").Copy() | Out-Null
    $ws.Range($ref).PasteSpecial(-4163) | Out-Null
}

# ---- Plain numeric value updates (style unchanged) ----
$ws.Range("I15").Value = 2
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = -33.333333333333
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 100
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -75
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = -50
$ws.Range("I16").Value = 13
$ws.Range("J16").Value = 17
$ws.Range("K16").Value = -23.529411764705
$ws.Range("L16").Value = -40.909090909090
$ws.Range("M16").Value = -27.777777777777
$ws.Range("N16").Value = -89.430894308943
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = 8.333333333333
$ws.Range("I17").Value = 33
$ws.Range("J17").Value = 28
$ws.Range("K17").Value = 17.857142857142
$ws.Range("L17").Value = 37.5
$ws.Range("M17").Value = 120
$ws.Range("N17").Value = 65
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = -14.285714285714
$ws.Range("I18").Value = 16
$ws.Range("J18").Value = 13
$ws.Range("K18").Value = 23.076923076923
$ws.Range("L18").Value = -40.740740740740
$ws.Range("M18").Value = 128.571428571429
$ws.Range("N18").Value = -71.428571428571
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 27
$ws.Range("G19").Value = 32
$ws.Range("H19").Value = -15.625
$ws.Range("I19").Value = 49
$ws.Range("J19").Value = 76
$ws.Range("K19").Value = -35.526315789473
$ws.Range("L19").Value = -57.017543859649
$ws.Range("M19").Value = 63.333333333333
$ws.Range("N19").Value = -10.909090909090
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 66.666666666666
$ws.Range("I20").Value = 9
$ws.Range("K20").Value = 28.571428571428
$ws.Range("L20").Value = 28.571428571428
$ws.Range("M20").Value = 12.5
$ws.Range("N20").Value = -83.928571428571
$ws.Range("C21").Value = 11
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = -54.166666666666
$ws.Range("F21").Value = 57
$ws.Range("G21").Value = 65
$ws.Range("H21").Value = -12.307692307692
$ws.Range("I21").Value = 122
$ws.Range("J21").Value = 143
$ws.Range("K21").Value = -14.685314685314
$ws.Range("L21").Value = -38.071065989847
$ws.Range("M21").Value = 50.617283950617
$ws.Range("N21").Value = -61.146496815286
$ws.Range("J22").Value = 2
$ws.Range("K22").Value = -50
$ws.Range("L22").Value = -83.333333333333
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = -25
$ws.Range("F23").Value = 8
$ws.Range("G23").Value = 15
$ws.Range("H23").Value = -46.666666666666
$ws.Range("I23").Value = 21
$ws.Range("J23").Value = 22
$ws.Range("K23").Value = -4.545454545454
$ws.Range("L23").Value = 23.529411764705
$ws.Range("M23").Value = -4.545454545454
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = 40
$ws.Range("F24").Value = 77
$ws.Range("G24").Value = 73
$ws.Range("H24").Value = 5.479452054794
$ws.Range("I24").Value = 165
$ws.Range("J24").Value = 164
$ws.Range("K24").Value = 0.609756097560
$ws.Range("L24").Value = -54.918032786885
$ws.Range("M24").Value = 71.875
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 20
$ws.Range("F25").Value = 28
$ws.Range("G25").Value = 29
$ws.Range("H25").Value = -3.448275862068
$ws.Range("I25").Value = 67
$ws.Range("J25").Value = 55
$ws.Range("K25").Value = 21.818181818181
$ws.Range("L25").Value = 39.583333333333
$ws.Range("M25").Value = 59.523809523809
$ws.Range("I26").Value = 4
$ws.Range("K26").Value = -20
$ws.Range("L26").Value = -42.857142857142
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 4
$ws.Range("J27").Value = 3
$ws.Range("K27").Value = 33.333333333333
$ws.Range("H30").Value = 0
$ws.Range("K30").Value = -50
$ws.Range("L30").Value = -66.666666666666

$excel.CutCopyMode = 0
